# multi grid actor path finder
# Adds a new "Size" (Int32) column (H) to the Actor sheet and two more
# data rows (Actor2 / Actor3), mirroring the existing Actor row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actor")

# --- New column H: header block (rows 2-4) ---------------------------------
$ws.Range("H2").Value = "All"
$ws.Range("H3").Value = "Size"
$ws.Range("H4").Value = "Int32"

# Row 4 (the "field type" row) uses a dedicated style (s="1", bold/alt font)
# on D4:G4 already - mirror that onto the new H4 cell by copying formats
# from the neighboring G4 cell.
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Existing data row 8 gains an H value -----------------------------------
$ws.Range("H8").Value = 1

# --- New data row 9 (Actor2) -------------------------------------------------
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "测试角色"
$ws.Range("C9").Value = "Warrior"
$ws.Range("D9").Value = 500
$ws.Range("E9").Value = 1100
$ws.Range("F9").Value = 1100
$ws.Range("G9").Value = "Exported/Actor/001/Prefabs/Actor2"
$ws.Range("H9").Value = 2

# --- New data row 10 (Actor3) ------------------------------------------------
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "测试角色"
$ws.Range("C10").Value = "Warrior"
$ws.Range("D10").Value = 500
$ws.Range("E10").Value = 1100
$ws.Range("F10").Value = 1100
$ws.Range("G10").Value = "Exported/Actor/001/Prefabs/Actor3"
$ws.Range("H10").Value = 3

# --- Selection moves to D14, matching the authored edit ---------------------
$ws.Range("D14").Select()
